# feat: add 2022-Q3 data
#
# 1. Duplicate the existing "2022-Q2" sheet (placed right after it) - the
#    duplicate keeps the old Q2 fund-holding data untouched and is renamed
#    back to "2022-Q2".
# 2. The original sheet is renamed to "2022-Q3" and its data is replaced
#    with the new Q3 fund holdings.
# 3. The "总计" (summary) sheet: row 2 is updated to the Q3 totals, and a
#    new row 3 is added holding the previous (Q2) totals.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# --- duplicate the "2022-Q2" sheet, placing the copy right after it ---
$wsQ2.Copy($null, $wsQ2)

# the original sheet (still referenced by $wsQ2) becomes "2022-Q3"
$wsQ2.Name = "2022-Q3"
$wsQ3 = $wsQ2

# the new copy (now in position 3) keeps the old data and becomes "2022-Q2"
$wb.Worksheets.Item(3).Name = "2022-Q2"

# =====================================================================
# update the "总计" summary sheet
# =====================================================================

# copy style from row 2 into new row 3 (same look as row 2 -> A column bold)
$wsTotal.Cells.Item(2,1).Copy($wsTotal.Cells.Item(3,1))

# row 3 gets the old 2022-Q2 totals (what used to be in row 2)
$wsTotal.Cells.Item(3,1).Value2 = 1
$wsTotal.Cells.Item(3,2).Value2 = "2022-Q2"
$wsTotal.Cells.Item(3,3).Value2 = 3
$wsTotal.Cells.Item(3,4).Value2 = 2

# row 2 is updated in place with the new 2022-Q3 totals
$wsTotal.Cells.Item(2,2).Value2 = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value2 = 3
$wsTotal.Cells.Item(2,4).Value2 = 1.97

# =====================================================================
# replace the data on the "2022-Q3" sheet with the new fund holdings
# =====================================================================

# wipe everything first (values + formatting) so we start clean
$wsQ3.Range("A1:H4").Clear()

# style used for the header row and the A column (bold, bordered, centered)
# -- same style already used on the "总计" sheet's header row
$styleSource = $wsTotal.Cells.Item(1,2)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $wsQ3.Cells.Item(1, $col)
    $styleSource.Copy($cell)
    $cell.Value2 = $headers[$col - 2]
}

$rowsQ3 = @(
    @(0, "001481", "华宝油气（QDII）美元",    "45.98", "94.53", "2.14", "0.9840", 5),
    @(1, "162411", "华宝油气（QDII）人民币A", "28.25", "94.53", "2.14", "0.6046", 5),
    @(2, "007844", "华宝油气（QDII）人民币 C","17.73", "94.53", "2.14", "0.3794", 5)
)

for ($i = 0; $i -lt $rowsQ3.Length; $i++) {
    $r = $i + 2
    $data = $rowsQ3[$i]

    # column A: numeric index, styled like the header (bold/border)
    $styleSource.Copy($wsQ3.Cells.Item($r,1))
    $wsQ3.Cells.Item($r,1).Value2 = $data[0]

    # columns B,C: plain text, default style
    $bRange = $wsQ3.Range($wsQ3.Cells.Item($r,2), $wsQ3.Cells.Item($r,3))
    $bRange.NumberFormat = "@"
    $wsQ3.Cells.Item($r,2).Value2 = $data[1]
    $wsQ3.Cells.Item($r,3).Value2 = $data[2]
    $bRange.ClearFormats()

    # columns D,E,F,G: numeric-looking values stored as text, default style
    $textRange = $wsQ3.Range($wsQ3.Cells.Item($r,4), $wsQ3.Cells.Item($r,7))
    $textRange.NumberFormat = "@"
    $wsQ3.Cells.Item($r,4).Value2 = $data[3]
    $wsQ3.Cells.Item($r,5).Value2 = $data[4]
    $wsQ3.Cells.Item($r,6).Value2 = $data[5]
    $wsQ3.Cells.Item($r,7).Value2 = $data[6]
    $textRange.ClearFormats()

    # column H: numeric rank, default style
    $wsQ3.Cells.Item($r,8).Value2 = $data[7]
}
